# ============================================================================
# Edit: "add note about meeting next week"
#
# 1. Title paragraph: drop the leading manual line-break run and the
#    now-stale _GoBack bookmark that wrapped "ALL IN".
# 2. Scripture text: cosmetic run-split (" " + "And when...") - no visible
#    text change, nothing to do.
# 3. Remove the "Close small group out with prayer requests." bullet.
# 4. "...registration window closes." -> "...registration window closes
#    next week. "
# 5. "...for a scholarship application." -> "...for a scholarship
#    application if needed."
# 6. Append a new highlighted paragraph announcing the leader meeting,
#    with the _GoBack bookmark now sitting inside it (between "plan" and
#    "s & expectations").
# ============================================================================

$d = $word.ActiveDocument

# --- 1. Title paragraph: remove manual line break + old _GoBack bookmark ---
$titlePara = $d.Paragraphs.Item(2)
$titleStart = $titlePara.Range.Start
$breakRange = $d.Range($titleStart, $titleStart + 1)
$breakRange.Delete()

$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# --- 3. Remove the "Close small group out with prayer requests." bullet ---
$d.Content.Find.Execute("Close small group out with prayer requests.`r", $true, $false, $false, $false, $false, $true, 1, $false, "", 1)

# --- 4. "...registration window closes." -> "...closes next week. " ------
$d.Content.Find.Execute("immediately before the registration window closes.", $true, $false, $false, $false, $false, $true, 1, $false, "immediately before the registration window closes next week. ", 2)

# --- 5. "...for a scholarship application." -> "...application if needed." -
$d.Content.Find.Execute("for a scholarship application.", $true, $false, $false, $false, $false, $true, 1, $false, "for a scholarship application if needed.", 2)

# --- 6. New paragraph: "Attention All Leaders: ..." -----------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$endRange = $lastPara.Range
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.ListFormat.RemoveNumbers()
$newPara.Style = "Normal"
$newPara.Range.ParagraphFormat.Reset()
$newPara.Range.Font.Name = "Arial"
$newPara.Range.Font.NameBi = "Arial"
$newPara.Range.Font.Size = 10
$newPara.Range.Font.SizeBi = 10

function Add-Seg($text, [bool]$bold) {
    $r = $word.ActiveDocument.Paragraphs.Item($word.ActiveDocument.Paragraphs.Count).Range
    $r.Collapse(0)
    $r.InsertAfter($text)
    $r.Font.Name = "Arial"
    $r.Font.NameBi = "Arial"
    $r.Font.Size = 10
    $r.Font.SizeBi = 10
    $r.Font.Bold = $bold
    $r.Font.BoldBi = $bold
    $r.HighlightColorIndex = 7
}

Add-Seg "Attention All Leaders" $true
Add-Seg ": there will be a MANDATORY Fusion Leader Meeting next Wednesday in the Chapel " $false
Add-Seg "at 7:00 PM" $false
Add-Seg " while students are in Large Group where we will discuss " $false
Add-Seg "Fusion plans & expectations" $false
Add-Seg " and answer all your questions." $false

# _GoBack bookmark now lives between "Fusion plan" and "s & expectations"
# (an empty / zero-length bookmark) - insert it last, via text search, so
# the position is never right at the very end of the document content.
$bkFind = $d.Content
$bkFind.Find.Execute("Fusion plans", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bkPos = $bkFind.End - 1
$bkRange = $d.Range($bkPos, $bkPos)
$d.Bookmarks.Add("_GoBack", $bkRange)

Write-Output "edit complete"
